# Add a new user record (row 5) to the users sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "849b9eb3-f482-40d0-9666-c0b0081c9534"
$ws.Range("B5").Value = "mustafa"
$ws.Range("C5").Value = "aka"
$ws.Range("D5").Value = "YWthMTIx"
